$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3, 4, 5 (keep header row 1 and AmazonSmokeTest row 2)
$ws.Range("A3:G5").EntireRow.Delete()

# Update row 2 values to the AmazonSmokeTest / BlockFile data
$ws.Range("A2").Value = "AmazonSmokeTest"
$ws.Range("B2").Value = "ON"
$ws.Range("C2").Value = "AmazonSmokeTest"
$ws.Range("D2").Value = "AmazonSmokeTest.xlsx"
$ws.Range("E2").Value = "AmazonObjects.properties"
$ws.Range("F2").Value = "AmazonBlocks.xlsx"
$ws.Range("G2").Value = "amazonLive"

# Update the selection / top-left cell to match the new active selection
$ws.Range("B1").Select()
